$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).ClearFormats()
}

$ws.Range("D2").Value = "53.657.73"
$ws.Range("E2").Value = "  -4.35%  "
$ws.Range("D3").Value = "2.223.48"
$ws.Range("E3").Value = "  -5.91%  "
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.22%  "
Set-TextValue "D5" "487.33"
$ws.Range("E5").Value = "  -3.02%  "
Set-TextValue "D6" "124.96"
$ws.Range("E6").Value = "  -3.42%  "
Set-TextValue "D7" "0.996"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  -4.26%  "
$ws.Range("D9").Value = "2.219.70"
$ws.Range("E9").Value = "  -6.19%  "
$ws.Range("E10").Value = "  -6.08%  "
$ws.Range("E11").Value = "  -1.28%  "
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("D14").Value = "2.613.49"
$ws.Range("E14").Value = "  -6.10%  "
Set-TextValue "D15" "21.13"
$ws.Range("E15").Value = "  -1.16%  "
$ws.Range("D16").Value = "53.600.88"
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("D18").Value = "2.212.83"
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D19" "9.57"
$ws.Range("E19").Value = "  -4.47%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D20" "3.96"
$ws.Range("E20").Value = "  -1.45%  "
Set-TextValue "D21" "295.02"
$ws.Range("E21").Value = "  -4.19%  "
Set-TextValue "D22" "6.17"
$ws.Range("E22").Value = "  -2.27%  "
Set-TextValue "D23" "0.999"
$ws.Range("E23").Value = "  -0.05%  "
Set-TextValue "D24" "62.84"
$ws.Range("E24").Value = "  -4.84%  "
Set-TextValue "D25" "0.998"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "2.321.64"
$ws.Range("E27").Value = "  -6.11%  "
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("E29").Value = "  -2.96%  "
Set-TextValue "D30" "165.16"
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("E31").Value = "  -3.86%  "
Set-TextValue "D32" "0.997"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E33").Value = "  -6.56%  "
Set-TextValue "D34" "0.996"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("E35").Value = "  -1.22%  "
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  -1.14%  "
Set-TextValue "D39" "0.838"
$ws.Range("E39").Value = "  +5.11%  "
Set-TextValue "D40" "3.55"
$ws.Range("E40").Value = "  -4.39%  "
Set-TextValue "D41" "35.83"
$ws.Range("E41").Value = "  -1.23%  "
Set-TextValue "D42" "0.368"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -1.12%  "
Set-TextValue "D44" "126.58"
$ws.Range("E44").Value = "  -2.01%  "
Set-TextValue "D45" "3.28"
$ws.Range("E45").Value = "  -2.50%  "
Set-TextValue "D46" "4.78"
$ws.Range("E46").Value = "  +2.01%  "
Set-TextValue "D47" "0.0881"
$ws.Range("E47").Value = "  -2.42%  "
Set-TextValue "D48" "0.534"
$ws.Range("E48").Value = "  -4.88%  "
Set-TextValue "D49" "231.77"
$ws.Range("E49").Value = "  -2.68%  "
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("E51").Value = "  -3.27%  "
